# Update MontosMotoboy_sushi_club price table with new values (columns B-H, rows 2-13)
# Column I is left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(4637, 4637, 4631, 4896, 4896, 4539, 4485)
    3  = @(4216, 4216, 4210, 4451, 4451, 4127, 4078)
    4  = @(4637, 4637, 4631, 4896, 4896, 4539, 4485)
    5  = @(4216, 4216, 4210, 4451, 4451, 4127, 4078)
    6  = @(5465, 5465, 5237, 5444, 5444, 6162, 5741)
    7  = @(5093, 5093, 4927, 5065, 5065, 5693, 5320)
    8  = @(4602, 4602, 4403, 4589, 4589, 5354, 4817)
    9  = @(4403, 4403, 4147, 4210, 4210, 5093, 4651)
    10 = @(6797, 6797, 6506, 6914, 6914, 7804, 7114)
    11 = @(6328, 6328, 6134, 6506, 6506, 7146, 6624)
    12 = @(5665, 5665, 5575, 5941, 5941, 6686, 5941)
    13 = @(5500, 5500, 5327, 5590, 5590, 6252, 5761)
}

foreach ($row in $newValues.Keys) {
    $cols = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        # Columns B through H correspond to index 2 through 8
        $colIndex = $i + 2
        $ws.Cells.Item($row, $colIndex).Value = $cols[$i]
    }
}

$wb.Save()
